# Auto-generated: apply scheduled market-data refresh to Titan Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 995
$ws.Range("I18").Value = 950
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 950
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -666
$ws.Range("N18").Value = -1568
$ws.Range("H19").Value = 1045.8823
$ws.Range("I19").Value = 889.0909
$ws.Range("K19").Value = 889.0909
$ws.Range("M19").Value = -714.0909
$ws.Range("H43").Value = 871.2727
$ws.Range("I43").Value = 850
$ws.Range("J43").Value = 883.4286
$ws.Range("K43").Value = 850
$ws.Range("L43").Value = 883.4286
$ws.Range("M43").Value = -781
$ws.Range("N43").Value = -1021.4286
$ws.Range("H107").Value = 463341.16
$ws.Range("I107").Value = 505372.2
$ws.Range("K107").Value = 505372.2
$ws.Range("M107").Value = -503452.2
$ws.Range("H113").Value = 43695.88
$ws.Range("I113").Value = 63447.35
$ws.Range("J113").Value = 1724
$ws.Range("K113").Value = 63447.35
$ws.Range("L113").Value = 1724
$ws.Range("M113").Value = -60193.35
$ws.Range("N113").Value = -8232
$ws.Range("H116").Value = 10646301
$ws.Range("I116").Value = 12581656
$ws.Range("J116").Value = 1850
$ws.Range("K116").Value = 12581656
$ws.Range("L116").Value = 1850
$ws.Range("M116").Value = -12578214
$ws.Range("N116").Value = -8734
$ws.Range("H129").Value = 1173.2858
$ws.Range("J129").Value = 1372.1875
$ws.Range("L129").Value = 4116.5625
$ws.Range("N129").Value = -14116.5625
$ws.Range("H132").Value = 214966.19
$ws.Range("I132").Value = 259243.94
$ws.Range("J132").Value = 54885.08
$ws.Range("K132").Value = 777731.8200000001
$ws.Range("L132").Value = 164655.24
$ws.Range("M132").Value = -775201.8200000001
$ws.Range("N132").Value = -169715.24
$ws.Range("H137").Value = 76924510
$ws.Range("I137").Value = 100001310
$ws.Range("K137").Value = 300003930
$ws.Range("M137").Value = -300001380

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2262.724
$ws.Range("I61").Value = 1743.25
$ws.Range("K61").Value = 1743.25
$ws.Range("M61").Value = -1531.25
$ws.Range("H74").Value = 11040.695
$ws.Range("I74").Value = 1646.8
$ws.Range("J74").Value = 73666.664
$ws.Range("K74").Value = 1646.8
$ws.Range("L74").Value = 73666.664
$ws.Range("M74").Value = -772.8
$ws.Range("N74").Value = -75414.664
$ws.Range("H77").Value = 11040.695
$ws.Range("I77").Value = 1646.8
$ws.Range("J77").Value = 73666.664
$ws.Range("K77").Value = 8234
$ws.Range("L77").Value = 368333.32
$ws.Range("M77").Value = -3866
$ws.Range("N77").Value = -377069.32
$ws.Range("H122").Value = 2417.0588
$ws.Range("I122").Value = 2227
$ws.Range("J122").Value = 2520.7273
$ws.Range("K122").Value = 6681
$ws.Range("L122").Value = 7562.1819
$ws.Range("M122").Value = -4231
$ws.Range("N122").Value = -12462.1819
$ws.Range("H136").Value = 2262.724
$ws.Range("I136").Value = 1743.25
$ws.Range("K136").Value = 5229.75
$ws.Range("M136").Value = -2679.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1376.1333
$ws.Range("J94").Value = 1236.75
$ws.Range("L94").Value = 1236.75
$ws.Range("N94").Value = -2138.75
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120
$ws.Range("H134").Value = 4856.4287
$ws.Range("I134").Value = 4260.5713
$ws.Range("J134").Value = 5154.357
$ws.Range("K134").Value = 12781.7139
$ws.Range("L134").Value = 15463.071
$ws.Range("M134").Value = -10246.7139
$ws.Range("N134").Value = -20533.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2777.75
$ws.Range("I16").Value = 555.5
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 555.5
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -268.5
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 1526.2106
$ws.Range("I31").Value = 884.6923
$ws.Range("J31").Value = 2916.1667
$ws.Range("K31").Value = 884.6923
$ws.Range("L31").Value = 2916.1667
$ws.Range("M31").Value = -589.6923
$ws.Range("N31").Value = -3506.1667
$ws.Range("H34").Value = 1526.2106
$ws.Range("I34").Value = 884.6923
$ws.Range("J34").Value = 2916.1667
$ws.Range("K34").Value = 884.6923
$ws.Range("L34").Value = 2916.1667
$ws.Range("M34").Value = -682.6923
$ws.Range("N34").Value = -3320.1667
$ws.Range("H58").Value = 940.4717000000001
$ws.Range("I58").Value = 559.70966
$ws.Range("J58").Value = 1477
$ws.Range("K58").Value = 559.70966
$ws.Range("L58").Value = 1477
$ws.Range("M58").Value = -356.70966
$ws.Range("N58").Value = -1883
$ws.Range("H99").Value = 5683084
$ws.Range("I99").Value = 6251242.5
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 6251242.5
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -6249744.5
$ws.Range("N99").Value = -4496
$ws.Range("H113").Value = 2777.75
$ws.Range("I113").Value = 555.5
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 555.5
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 1614.5
$ws.Range("N113").Value = -9340
$ws.Range("H126").Value = 5683084
$ws.Range("I126").Value = 6251242.5
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 18753727.5
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -18751257.5
$ws.Range("N126").Value = -9440
$ws.Range("H132").Value = 1817.7667
$ws.Range("I132").Value = 1479.1555
$ws.Range("K132").Value = 4437.4665
$ws.Range("M132").Value = -1907.4665
$ws.Range("H134").Value = 3217.5625
$ws.Range("I134").Value = 1777.75
$ws.Range("J134").Value = 5617.25
$ws.Range("K134").Value = 5333.25
$ws.Range("L134").Value = 16851.75
$ws.Range("M134").Value = -2798.25
$ws.Range("N134").Value = -21921.75
$ws.Range("H136").Value = 940.4717000000001
$ws.Range("I136").Value = 559.70966
$ws.Range("J136").Value = 1477
$ws.Range("K136").Value = 1679.12898
$ws.Range("L136").Value = 4431
$ws.Range("M136").Value = 870.87102
$ws.Range("N136").Value = -9531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1632.4584
$ws.Range("J131").Value = 1807.0952
$ws.Range("L131").Value = 5421.2856
$ws.Range("N131").Value = -15501.2856
$ws.Range("H137").Value = 7773161
$ws.Range("I137").Value = 10003387
$ws.Range("J137").Value = 339074.34
$ws.Range("K137").Value = 30010161
$ws.Range("L137").Value = 1017223.02
$ws.Range("M137").Value = -30005061
$ws.Range("N137").Value = -1027423.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H70").Value = 6656.4116
$ws.Range("I70").Value = 7864.3335
$ws.Range("J70").Value = 5297.5
$ws.Range("K70").Value = 7864.3335
$ws.Range("L70").Value = 5297.5
$ws.Range("M70").Value = -7594.3335
$ws.Range("N70").Value = -5837.5
$ws.Range("H73").Value = 6656.4116
$ws.Range("I73").Value = 7864.3335
$ws.Range("J73").Value = 5297.5
$ws.Range("K73").Value = 7864.3335
$ws.Range("L73").Value = 5297.5
$ws.Range("M73").Value = -6928.3335
$ws.Range("N73").Value = -7169.5
$ws.Range("H132").Value = 5079.9585
$ws.Range("I132").Value = 5941.2
$ws.Range("J132").Value = 4464.7856
$ws.Range("K132").Value = 17823.6
$ws.Range("L132").Value = 13394.3568
$ws.Range("M132").Value = -15293.6
$ws.Range("N132").Value = -18454.3568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2728.8064
$ws.Range("I7").Value = 1591.1818
$ws.Range("J7").Value = 3354.5
$ws.Range("K7").Value = 1591.1818
$ws.Range("L7").Value = 3354.5
$ws.Range("M7").Value = -1479.1818
$ws.Range("N7").Value = -3578.5
$ws.Range("H40").Value = 3002.4546
$ws.Range("I40").Value = 1300.2142
$ws.Range("K40").Value = 1300.2142
$ws.Range("M40").Value = -1164.2142
$ws.Range("H122").Value = 3220.9285
$ws.Range("I122").Value = 2144.5557
$ws.Range("J122").Value = 3730.7896
$ws.Range("K122").Value = 6433.6671
$ws.Range("L122").Value = 11192.3688
$ws.Range("M122").Value = -3983.6671
$ws.Range("N122").Value = -16092.3688
$ws.Range("H126").Value = 2728.8064
$ws.Range("I126").Value = 1591.1818
$ws.Range("J126").Value = 3354.5
$ws.Range("K126").Value = 4773.5454
$ws.Range("L126").Value = 10063.5
$ws.Range("M126").Value = -2303.5454
$ws.Range("N126").Value = -15003.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 55255.527
$ws.Range("I122").Value = 60873.824
$ws.Range("K122").Value = 182621.472
$ws.Range("M122").Value = -180171.472
$ws.Range("H126").Value = 74579.14
$ws.Range("I126").Value = 86400.664
$ws.Range("K126").Value = 259201.992
$ws.Range("M126").Value = -256731.992
$ws.Range("H132").Value = 11365508
$ws.Range("I132").Value = 18520308
$ws.Range("J132").Value = 2004.1177
$ws.Range("K132").Value = 55560924
$ws.Range("L132").Value = 6012.3531
$ws.Range("M132").Value = -55558394
$ws.Range("N132").Value = -11072.3531

